# Updated symbol list on Tue Dec 20 21:23:26 UTC 2022 with GitHub Actions
#
# Refreshes the "Price" (column D) values for the listed coins, and tweaks
# two "Volume(1h)" (column E) labels. All of these cells hold plain text in
# the workbook (e.g. "--" appears for missing prices), so the numeric-looking
# price strings must be written back as text rather than being allowed to
# convert into numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row number -> new Price (column D) text
$priceUpdates = @{
    2  = "250.12"
    3  = "22.91"
    4  = "5.451"
    5  = "0.05660"
    6  = "3.416"
    7  = "6.383"
    8  = "0.8152"
    9  = "0.9304"
    10 = "0.1441"
    11 = "0.07503"
    13 = "0.03092"
    14 = "0.09360"
    15 = "3.759"
    16 = "0.001591"
    18 = "0.0005797"
    19 = "0.006411"
    20 = "0.005038"
    21 = "0.001033"
    23 = "3.702"
    24 = "2.178"
    25 = "0.3302"
    26 = "0.1282"
    28 = "0.0003004"
    40 = "0.04019"
    41 = "0.006783"
    42 = "0.1068"
    43 = "0.002711"
    44 = "0.008033"
    45 = "0.00005810"
    46 = "0.00000000750"
    47 = "0.5006"
    49 = "0.00002101"
}

foreach ($row in $priceUpdates.Keys) {
    $cell = $ws.Range("D$row")
    # Force text formatting so the numeric-looking string is not coerced
    # into a Number value by Excel.
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$row]
}

# row number -> new Volume(1h) (column E) text
$volumeUpdates = @{
    18 = "17OneONEWorstin24h"
    47 = "46CoinbaseStockTokenCOIN"
}

foreach ($row in $volumeUpdates.Keys) {
    $ws.Range("E$row").Value = $volumeUpdates[$row]
}
